{"js": "// Update each math-fact cell and the header date in-place, preserving\n// run formatting, by searching for the exact old text and replacing it.\nconst replacements = [\n  [\"2026-01-14 Wednesday\", \"2026-01-15 Thursday\"],\n  [\"54\u00d755=2970\", \"26\u00d745=1170\"],\n  [\"28\u00d798=2744\", \"57\u00d734=1938\"],\n  [\"34\u00d726=884\", \"69\u00d737=2553\"],\n  [\"61\u00d735=2135\", \"53\u00d788=4664\"],\n  [\"88\u00d775=6600\", \"36\u00d756=2016\"],\n  [\"70\u00d743=3010\", \"31\u00d796=2976\"],\n  [\"53\u00d797=5141\", \"46\u00d785=3910\"],\n  [\"95\u00d753=5035\", \"92\u00d754=4968\"],\n  [\"52\u00d734=1768\", \"84\u00d741=3444\"],\n  [\"81\u00d776=6156\", \"96\u00d795=9120\"],\n  [\"63\u00d748=3024\", \"46\u00d723=1058\"],\n  [\"97\u00d788=8536\", \"48\u00d721=1008\"],\n  [\"60\u00d753=3180\", \"30\u00d749=1470\"],\n  [\"79\u00d776=6004\", \"62\u00d770=4340\"],\n  [\"19\u00d736=684\", \"80\u00d750=4000\"],\n  [\"12\u00d762=744\", \"22\u00d786=1892\"],\n  [\"72\u00d717=1224\", \"37\u00d767=2479\"],\n  [\"11\u00d785=935\", \"70\u00d791=6370\"],\n  [\"23\u00d768=1564\", \"24\u00d793=2232\"],\n  [\"49\u00d779=3871\", \"96\u00d799=9504\"],\n  [\"95\u00d774=7030\", \"18\u00d755=990\"],\n  [\"57\u00d794=5358\", \"59\u00d789=5251\"],\n  [\"38\u00d758=2204\", \"42\u00d751=2142\"],\n  [\"92\u00d794=8648\", \"12\u00d773=876\"],\n  [\"82\u00d725=2050\", \"15\u00d741=615\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and each two-digit-by-two-digit multiplication\n# fact in the table by running Find/Replace over the whole document for\n# each exact old value, preserving the surrounding run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2026-01-14 Wednesday\"; New = \"2026-01-15 Thursday\" }\n    @{ Old = \"54\u00d755=2970\"; New = \"26\u00d745=1170\" }\n    @{ Old = \"28\u00d798=2744\"; New = \"57\u00d734=1938\" }\n    @{ Old = \"34\u00d726=884\"; New = \"69\u00d737=2553\" }\n    @{ Old = \"61\u00d735=2135\"; New = \"53\u00d788=4664\" }\n    @{ Old = \"88\u00d775=6600\"; New = \"36\u00d756=2016\" }\n    @{ Old = \"70\u00d743=3010\"; New = \"31\u00d796=2976\" }\n    @{ Old = \"53\u00d797=5141\"; New = \"46\u00d785=3910\" }\n    @{ Old = \"95\u00d753=5035\"; New = \"92\u00d754=4968\" }\n    @{ Old = \"52\u00d734=1768\"; New = \"84\u00d741=3444\" }\n    @{ Old = \"81\u00d776=6156\"; New = \"96\u00d795=9120\" }\n    @{ Old = \"63\u00d748=3024\"; New = \"46\u00d723=1058\" }\n    @{ Old = \"97\u00d788=8536\"; New = \"48\u00d721=1008\" }\n    @{ Old = \"60\u00d753=3180\"; New = \"30\u00d749=1470\" }\n    @{ Old = \"79\u00d776=6004\"; New = \"62\u00d770=4340\" }\n    @{ Old = \"19\u00d736=684\"; New = \"80\u00d750=4000\" }\n    @{ Old = \"12\u00d762=744\"; New = \"22\u00d786=1892\" }\n    @{ Old = \"72\u00d717=1224\"; New = \"37\u00d767=2479\" }\n    @{ Old = \"11\u00d785=935\"; New = \"70\u00d791=6370\" }\n    @{ Old = \"23\u00d768=1564\"; New = \"24\u00d793=2232\" }\n    @{ Old = \"49\u00d779=3871\"; New = \"96\u00d799=9504\" }\n    @{ Old = \"95\u00d774=7030\"; New = \"18\u00d755=990\" }\n    @{ Old = \"57\u00d794=5358\"; New = \"59\u00d789=5251\" }\n    @{ Old = \"38\u00d758=2204\"; New = \"42\u00d751=2142\" }\n    @{ Old = \"92\u00d794=8648\"; New = \"12\u00d773=876\" }\n    @{ Old = \"82\u00d725=2050\"; New = \"15\u00d741=615\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
